$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B152").Value = "[name=""フロストノヴァ""]  一人、撃退した……！ だが、まさか……奴は自ら撤退したのか？`n"
$ws.Range("C3").Value = "[name=""'Emperors' Blade'""]  Such potential!`n"
$ws.Range("C4").Value = "[name=""'Emperors' Blade'""]  *Hiss*...! I should take my hat off to you! `n"
$ws.Range("C5").Value = "[name=""'Emperors' Blade'""]  Alas, I forgot to bring my cap with me today.`n"
$ws.Range("C7").Value = "[name=""'Emperors' Blade'""]  It's pointless to keep fighting.`n"
$ws.Range("C8").Value = "[name=""'Emperors' Blade'""]  We are merely protectors of Ursus. You may take you men with you. Live on as the protector of these Infected if that is your wish. Let us take what we must to realize our respective goals. I find that reasonable.   `n"
$ws.Range("C10").Value = "[name=""'Emperors' Blade'""]  All nations have men who call themselves protectors. There are many things that need to be protected, and also far too many men who call themselves protectors. `n"
$ws.Range("C11").Value = "[name=""'Emperors' Blade'""]  Most of them are corrupt and incompetent, yet the rest of us excel precisely because of the things we protect.`n"
$ws.Range("C12").Value = "[name=""'Emperors' Blade'""]  We are protecting Ursus's future.`n"
$ws.Range("C17").Value = "[name=""'Emperors' Blade'""]  ......`n"
$ws.Range("C48").Value = "[name=""'Emperors' Blade'""]  ......`n"
$ws.Range("C77").Value = "[name=""'Emperors' Blade'""]  ......`n"
$ws.Range("C18").Value = "[name=""'Emperors' Blade'""]  All societies will at some point have to bring out their lashes and whip their laziest and most simpleminded members. That being said, we are blades, not whips. `n"
$ws.Range("C19").Value = "[name=""'Emperors' Blade'""]  If cutting off branches that do nothing but suck the tree dry of nutrients was a crime, then yes, that makes us evil.`n"
$ws.Range("C20").Value = "[name=""'Emperors' Blade'""]  You said you will 'expose' our 'atrocities'. That merely tells me you are unable to confront us, not to mention how likely the way you see things will change.  `n"
$ws.Range("C21").Value = "[name=""'Emperors' Blade'""]  If I walk towards you, I suppose you will take a step back?`n"
$ws.Range("C24").Value = "[name=""'Emperors' Blade'""]  Ursus's beneficence is of equal weight to its devilries. When you stand up against our evil deeds, you must also face all the good that we have done.`n"
$ws.Range("C25").Value = "[name=""'Emperors' Blade'""]  All of it.`n"
$ws.Range("C26").Value = "[name=""'Emperors' Blade'""]  A nation cannot be measured by good or evil. Nothing is worthy of discussion on this metric. `n"
$ws.Range("C27").Value = "[name=""'Emperors' Blade'""]  ...Perhaps it's not realistic to expect you to reach his heights and fully grasp his wisdom so soon.`n"
$ws.Range("C28").Value = "[name=""'Emperors' Blade'""]  But if everything is as he said, then perhaps we will come for your advice when that day comes.`n"
$ws.Range("C29").Value = "[name=""'Emperors' Blade'""]  Right now, aside the potential and the seeds hidden in you, you have nothing.`n"
$ws.Range("C30").Value = "[name=""'Emperors' Blade'""]  The empire isn't something you can understand, Vouivre.`n"
$ws.Range("C31").Value = "[name=""'Emperors' Blade'""]  Wait. Vouivre... No... `n"
$ws.Range("C32").Value = "[name=""'Emperors' Blade'""]  You... It can't be...`n"
$ws.Range("C33").Value = "[name=""'Emperors' Blade'""]  36... 【Code word】!`n"
$ws.Range("C34").Value = "[name=""'Emperors' Blade'""]  —Wendigo—`n"
$ws.Range("C38").Value = "[name=""'Emperors' Blade'""]  No. Patriot. No.  `n"
$ws.Range("C41").Value = "[name=""'Emperors' Blade'""]  We don't want to be your enemy!`n"
$ws.Range("C42").Value = "[name=""'Emperors' Blade'""]  Wendigo... You are a legend known to a few of us in the Empire's army.`n"
$ws.Range("C43").Value = "[name=""'Emperors' Blade'""]  Even if the common folk of the nomadic cities have long forgotten, all of us still remember the tales our forebears told us.`n"
$ws.Range("C44").Value = "[name=""'Emperors' Blade'""]  You have my respect, Wendigo!`n"
$ws.Range("C47").Value = "[name=""'Emperors' Blade'""]  ...So you follow her. Lying has become second nature to the soldiers of the tundras, it's hard to know what to believe. `n"
$ws.Range("C49").Value = "[name=""'Emperors' Blade'""]  So it's the truth, then. You are Infected.`n"
$ws.Range("C51").Value = "[name=""'Emperors' Blade'""]  That's not right, Wendigo. This plan will fail. `n"
$ws.Range("C53").Value = "[name=""'Emperors' Blade'""]  Fantasies do not come true just because more people believe in them.`n"
$ws.Range("C55").Value = "[name=""'Emperors' Blade'""]  You're right. My generation hasn't seen much action. `n"
$ws.Range("C56").Value = "[name=""'Emperors' Blade'""]  Therefore... I am inviting you to join us. Ursus Captain Buldrokkas'tee, I ask you to come with us. Ursus needs you.  `n"
$ws.Range("C59").Value = "[name=""'Emperors' Blade'""]  Having fought under the command of the last emperor for over a hundred years, surely you must remember Ursus's greatness in those times.`n"
$ws.Range("C60").Value = "[name=""'Emperors' Blade'""]  Was it a prosperous and marvelous era?`n"
$ws.Range("C61").Value = "[name=""'Emperors' Blade'""]  Everyone united under the banner of Ursus regardless of race, and they all fought for its future. Enemies far and wide fell to our blades and cannon fire, and we stood undefeated in the battlefield.`n"
$ws.Range("C62").Value = "[name=""'Emperors' Blade'""]  I am far from the only one who wishes to return to those glorious times... All of us want to see again an era when all men are brothers and share the same deep hatred for their enemies.`n"
$ws.Range("C63").Value = "[name=""'Emperors' Blade'""]  We can bring all those greedy countries to part with the lands they stole, and give the downtrodden the chance to regain their dignity under the glory of Ursus. We conquer not to annihilate, but to rebuild.`n"
$ws.Range("C64").Value = "[name=""'Emperors' Blade'""]  We gave this world a new lease of life.`n"
$ws.Range("C65").Value = "[name=""'Emperors' Blade'""]  Let us rebuild that era. With all people united, facing all the storms that come our way together. I do believe it's far better than living your days on the tundra, toiling away. `n"
$ws.Range("C67").Value = "[name=""'Emperors' Blade'""]  No one has it easy. Such is the cruelty of the times we live in. Life is difficult for Ursus citizens as well.  `n"
$ws.Range("C68").Value = "[name=""'Emperors' Blade'""]  The lack of order, the loss of power, and the destruction of morals. These are all mistakes.`n"
$ws.Range("C69").Value = "[name=""'Emperors' Blade'""]  These are the mistakes that ruined Ursus, and we know where the problem lies.`n"
$ws.Range("C70").Value = "[name=""'Emperors' Blade'""]  We can fix these problems.`n"
$ws.Range("C72").Value = "[name=""'Emperors' Blade'""]  You may ask your leader and see if he believes us.`n"
$ws.Range("C74").Value = "[name=""'Emperors' Blade'""]  Buldrokkas'tee, we can set Ursus on the right path.`n"
$ws.Range("C80").Value = "[name=""'Emperors' Blade'""]  You cannot deny the glory of that era. Everything you did tied you back to Ursus. Even your title symbolizes your desires and aspirations. `n"
$ws.Range("C83").Value = "[name=""'Emperors' Blade'""]  Talulah... *Hiss*.`n"
$ws.Range("C84").Value = "[name=""'Emperors' Blade'""]  ...You are right.`n"
$ws.Range("C85").Value = "[name=""'Emperors' Blade'""]  Perhaps things will be different with the help of you and your men.`n"
$ws.Range("C87").Value = "[name=""'Emperors' Blade'""]  Blades do not make promises. Weapons do not ever make promises. `n"
$ws.Range("C88").Value = "[name=""'Emperors' Blade'""]  But I believe the Infected are a force that Ursus should've possessed in the first place.`n"
$ws.Range("C89").Value = "[name=""'Emperors' Blade'""]  All of you should strive to become a symbol of Ursus's glory. The Infected who oppose Ursus shall perish. As for all of you, if you are willing to take up arms for Ursus, it will be recognized as a glorious and noble deed.  `n"
$ws.Range("C93").Value = "[name=""'Emperors' Blade'""]  In that case, if you are fighting for the minority, on what basis are you saying you deserve the majority's endorsement? What is so righteous about your fight for the Infected?  `n"
$ws.Range("C94").Value = "[name=""'Emperors' Blade'""]  Surely you will be able to give me an intelligent answer, Wendigo.`n"
$ws.Range("C97").Value = "[name=""'Emperors' Blade'""]  Hfff...`n"
$ws.Range("C99").Value = "[name=""'Emperors' Blade'""]  If I answered, 'It had nothing to do with us...'`n"
$ws.Range("C100").Value = "[name=""'Emperors' Blade'""]  ...Hfff...`n"
$ws.Range("C101").Value = "[name=""'Emperors' Blade'""]  It looks like we won't be coming to an understanding.`n"
$ws.Range("C105").Value = "[name=""'Emperors' Blade'""]  *Hiss*...`n"
$ws.Range("C106").Value = "[name=""'Emperors' Blade'""]  Why her?`n"
$ws.Range("C110").Value = "[name=""'Emperors' Blade'""]  We, too, need the support of the Infected, and we are also fighting for a just cause.`n"
$ws.Range("C111").Value = "[name=""'Emperors' Blade'""]  Why do you refuse us, then?`n"
$ws.Range("C112").Value = "[name=""'Emperors' Blade'""]  Is it because you believe the path the 'Reunion' is taking you is better than the unity and reform we seek? Because it's more efficient? Or is it perhaps because it hurts the Ursus common folk less?  `n"
$ws.Range("C116").Value = "[name=""'Emperors' Blade'""]  It is our hope the Infected will unite under your lead.`n"
$ws.Range("C119").Value = "[name=""Patriot""]  I've already walked the path you are singing praises of, and that's why I don't believe there are any 'better choices' out there. If you tell me you had that kind of foresight, it just tells me you haven't yet faced the irony of fate. `n"
$ws.Range("C122").Value = "[name=""'Emperors' Blade'""]  A declaration is nothing more than that, a declaration.`n"
$ws.Range("C123").Value = "[name=""'Emperors' Blade'""]  Even if that is what you believe, Wendigo, what about the people standing behind you? Will they trust her as much as you do? `n"
$ws.Range("C124").Value = "[name=""'Emperors' Blade'""]  Even if they revere your might and your honesty, how will you assure them they are not just worshipping a powerful idol?`n"
$ws.Range("C125").Value = "[name=""'Emperors' Blade'""]  Is an evil idol able to shelter them better, or is an idol who lets herself be trampled upon more preferable?`n"
$ws.Range("C127").Value = "[name=""'Emperors' Blade'""]  Reality is far colder than the flurries of the Northwestern Tundra, Wendigo. Once you lose power... you'll find yourself at the mercy of those who are capable. And the first ones to make that move might not be your enemies.  `n"
$ws.Range("C128").Value = "[name=""'Emperors' Blade'""]  Perhaps they know not who Kashchey is, but they know what a duke is, and they know the duke's daughter will one day become the duchess.`n"
$ws.Range("C129").Value = "[name=""'Emperors' Blade'""]  The people around you are not like how you see them. Her strength will one day come to rival yours, she will be no less resourceful, and her cunning will be on another level. `n"
$ws.Range("C131").Value = "[name=""'Emperors' Blade'""]  No, Wendigo, no... You haven't been an Ursus citizen long enough. You don't understand. I am talking about Kashchey.  `n"
$ws.Range("C133").Value = "[name=""'Emperors' Blade'""]  *Hiss*... *Huff*...`n"
$ws.Range("C134").Value = "[name=""'Emperors' Blade'""]  ...I've heard stories that you dote on young people, or perhaps pamper is the better term.`n"
$ws.Range("C137").Value = "[name=""'Emperors' Blade'""]  Buldrokkas'tee...`n"
$ws.Range("C141").Value = "[name=""'Emperors' Blade'""]  ...Ah. Wendigo... I'm afraid I cannot spill any more secrets. `n"
$ws.Range("C142").Value = "[name=""'Emperors' Blade'""]  Our warning isn't merely a rhetorical tool, Buldrokkas'tee.`n"
$ws.Range("C143").Value = "[name=""'Emperors' Blade'""]  We can only mourn for you.`n"
$ws.Range("C145").Value = "[name=""'Emperors' Blade'""]  That day won't come, Wendigo. That day won't come. `n"
$ws.Range("C146").Value = "[name=""'Emperors' Blade'""]  Farewell, 'Patriot'. This day has been unfruitful. How regretful. `n"
$ws.Range("C147").Value = "[name=""'Emperors' Blade'""]  I should add, Talulah... We look forward to seeing your improvement on the battlefield.`n"
$ws.Range("C148").Value = "[name=""'Emperors' Blade'""]  You should lead them. You are more than up to the task. `n"
$ws.Range("D30").Value = "[name="" '황제의 칼날'""]  제국은 네가 이해할 수 없는 것이다, 뷔브르.`n"
$ws.Range("D31").Value = "[name="" '황제의 칼날'""]  잠깐. 뷔브르…… 아니……`n"
